# Refresh the cryptocurrency price/volume snapshot (GitHub Actions nightly update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the "Price" (column D) and "Volume(1h)" (column E) cells,
# keyed by cell address, in row order.
$newValues = [ordered]@{
    "D2" = "27.282.88"
    "E2" = "  +1.07%  "
    "D3" = "1.566.44"
    "E3" = "  +0.16%  "
    "E5" = "  +1.33%  "
    "E6" = "  +0.12%  "
    "E7" = "  -0.28%  "
    "D8" = "22.05"
    "E8" = "  -0.26%  "
    "E9" = "  +0.30%  "
    "E10" = "  -0.48%  "
    "E11" = "  +1.73%  "
    "D12" = "1.789.70"
    "E12" = "  +0.17%  "
    "D13" = "1.556.53"
    "E13" = "  -0.43%  "
    "E14" = "  +0.56%  "
    "D16" = "27.267.27"
    "E16" = "  +1.02%  "
    "D17" = "61.90"
    "E17" = "  -0.18%  "
    "D18" = "217.19"
    "E18" = "  +0.54%  "
    "E19" = "  +1.54%  "
    "E20" = "  -0.10%  "
    "E21" = "  -0.28%  "
    "D22" = "4.13"
    "E22" = "  +0.61%  "
    "D23" = "9.23"
    "E23" = "  +0.35%  "
    "E24" = "  +0.08%  "
    "D25" = "153.11"
    "E25" = "  +0.07%  "
    "E26" = "  +0.29%  "
    "D27" = "15.00"
    "E27" = "  -0.60%  "
    "D28" = "0.106"
    "E28" = "  +1.55%  "
    "E29" = "  -0.15%  "
    "E30" = "  +2.03%  "
    "D31" = "0.0471"
    "E31" = "  +0.04%  "
    "E32" = "  +0.20%  "
    "D33" = "3.15"
    "E33" = "  +1.12%  "
    "D34" = "1.436.03"
    "E34" = "  +0.89%  "
    "E35" = "  +3.48%  "
    "E36" = "  +0.19%  "
    "E37" = "  -0.19%  "
    "E39" = "  +0.23%  "
    "D40" = "5.94"
    "E40" = "  +2.10%  "
    "D41" = "0.807"
    "E41" = "  +0.02%  "
    "E42" = "  -0.23%  "
    "D43" = "2.33"
    "E43" = "  +0.45%  "
    "D44" = "0.998"
    "E44" = "  -0.73%  "
    "D45" = "64.52"
    "E45" = "  -0.36%  "
    "E46" = "  -0.84%  "
    "D47" = "1.703.16"
    "E47" = "  +0.18%  "
    "D48" = "86.05"
    "E48" = "  -1.43%  "
    "E49" = "  +1.11%  "
    "E50" = "  +1.51%  "
    "D51" = "0.0954"
    "E51" = "  -0.54%  "
}

# A handful of "Price" cells hold plain decimals (e.g. "22.08", "0.998").
# Excel auto-converts a numeric-looking string typed into a General-format
# cell into a real number, which would silently drop significant trailing
# zeros (e.g. "15.00" -> 15). Mark just those cells as Text first so the
# exact string is preserved, matching how the source data is published.
$forceTextCells = @(
    "D8",
    "D17",
    "D18",
    "D22",
    "D23",
    "D25",
    "D27",
    "D28",
    "D31",
    "D33",
    "D40",
    "D41",
    "D43",
    "D44",
    "D45",
    "D48",
    "D51"
)
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

# Restore the default ("Normal") style on the cells that were temporarily
# switched to Text format, so no visible formatting change is left behind.
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
